# The presentation's embedded theme parts (ppt/theme/theme1.xml and
# ppt/theme/theme2.xml) had their colour schemes swapped: theme1.xml
# (used only by the Notes Master) becomes the "Integral" / Red Violet
# scheme, and theme2.xml (used by the Slide Master + the presentation
# itself) becomes the plain "Office" colour scheme.
#
# The object model only exposes one editable palette — the active
# design's ThemeColorScheme, which is backed by the Slide Master's
# theme part — so we push the target "Office" palette into it here.

function HexToRgbVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# Order of ThemeColorScheme.Item(n): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = HexToRgbVal $officeColors[$i - 1]
}
